$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AR" (AR model calibration)
# ---------------------------------------------------------------------------
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B2").Value = -0.0155998373882116
$wsAR.Range("B3").Value = 0.8080693289591837
$wsAR.Range("B4").Value = 49.36903627112625
$wsAR.Range("B5").Value = "[1.0, 0.19023100517363417, 0.15957736476760695, 0.16915007635191742, 0.13025529307705677, 0.3651582042972734, 0.1375205118505719, 0.15144446280071494, 0.1574439255637781, 0.09053995764706509, 0.09281341298415349, 0.09823517080724811, 0.0907387799556728, 0.1194926198492851, 0.1416753636656969, 0.1054058682648264, 0.10348847331058442, 0.0636944224454599, 0.08466370979157527, 0.07991439927286213]"

# ---------------------------------------------------------------------------
# Sheet "SETAR" (SETAR model calibration)
# ---------------------------------------------------------------------------
$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B3").Value = 0.9108910891089109
$wsSETAR.Range("B4").Value = -2.660583875057185
$wsSETAR.Range("B5").Value = 0.7124804935154867
$wsSETAR.Range("B6").Value = 36.3794048605352
$wsSETAR.Range("B7").Value = 2.554970230688115
$wsSETAR.Range("B8").Value = 0.6909033534613436
$wsSETAR.Range("B9").Value = 30.43277575694534
$wsSETAR.Range("B10").Value = "[0.9999999999999999, 0.2117185060040424, 0.1868486722447619, 0.185671437039581, 0.19014342083429547, 0.24196770189726288, 0.1146872154050553, 0.15148431535082427, 0.15019020761816368, 0.11027954481969852, 0.1299340423183486, 0.11067989939237509, 0.12248137163331348, 0.13327679548900925, 0.14627548145225439, 0.08708399970625862, 0.0981745480129427, 0.09602965718752753, 0.10184818520682139, 0.06962746290834544]"

# ---------------------------------------------------------------------------
# Sheet "GARCH" (GARCH model calibration)
# ---------------------------------------------------------------------------
$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B2").Value = -0.03977126274481312
$wsGARCH.Range("B3").Value = 1.311842050611887
$wsGARCH.Range("B4").Value = 0.08015125752694374
$wsGARCH.Range("B5").Value = 0.8963444262047099
$wsGARCH.Range("B6").Value = "[1.0, -0.010808380571508045, -0.02674827873927646, -0.024362968918563807, -0.05807573045685614, 0.30877310904565675, -0.02069091170190842, 0.014038259783707664, 0.008003352373506505, -0.040323179996447744, -0.04700525090030734, -0.05183819264679333, -0.04341594451254616, -0.011389003153922252, 0.032540767691482626, -0.017687163600546483, -0.008399396138445661, -0.04869495253717731, -0.04715054818403407, -0.013961750501083288]"

# ---------------------------------------------------------------------------
# Sheet "TARCH" (TARCH model calibration)
# ---------------------------------------------------------------------------
$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B2").Value = -0.06199952145114571
$wsTARCH.Range("B3").Value = 1.310222149942387
$wsTARCH.Range("B4").Value = 0.07416584059307202
$wsTARCH.Range("B5").Value = 0.896916432021676
$wsTARCH.Range("B6").Value = "[1.0, -0.010661196559691874, -0.026726530615183108, -0.024180422023687226, -0.05766887016402388, 0.3086060570930207, -0.02029259285375295, 0.014300440313646176, 0.008322663918891287, -0.03981838353263321, -0.046521164121585704, -0.051670852733149406, -0.04316146860662575, -0.01136178792940792, 0.03263703658557904, -0.017615061700667328, -0.008336309685379933, -0.04846878248876332, -0.04690195997255957, -0.013956762647933203]"
$wsTARCH.Range("B7").Value = 0.01072889588444607

# ---------------------------------------------------------------------------
# Sheet "AR_TARCH" (AR_TARCH model calibration)
# ---------------------------------------------------------------------------
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")
$wsARTARCH.Range("B2").Value = -0.02475887293644817
$wsARTARCH.Range("B3").Value = 1.285085415909744
$wsARTARCH.Range("B4").Value = 0.08816511381665533
$wsARTARCH.Range("B5").Value = 0.8903285959973737
$wsARTARCH.Range("B6").Value = "[1.0, 0.01202961087276561, 0.0011238303180410277, -0.0007174629323332885, -0.0303059370476774, 0.23386648165854698, -0.022241612781317633, 0.007632277479983158, -0.00023023096752381806, -0.03604941000908999, -0.05063782453733024, -0.03987532817805389, -0.031924615828510264, -0.018250414152227837, 0.027547719324608352, -0.010198708742010551, -0.0017548935395702908, -0.0344768342298025, -0.04153134719806645, -0.029408173527027487]"
$wsARTARCH.Range("B7").Value = -0.008449193637179889
$wsARTARCH.Range("B9").Value = 0.7988935808046292

Write-Output "Applied calibration updates across AR, SETAR, GARCH, TARCH, AR_TARCH sheets"
